$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

function Find-ParagraphIndex($exactText) {
    $want = $exactText + [char]13
    for ($k = 1; $k -le $d.Paragraphs.Count; $k++) {
        if ($d.Paragraphs.Item($k).Range.Text -eq $want) {
            return $k
        }
    }
    return -1
}

# Inserts a brand-new standalone paragraph, with text $newText, immediately
# before the (existing, untouched) paragraph whose exact text is $beforeText.
function Insert-ParagraphBefore($beforeText, $newText) {
    $i = Find-ParagraphIndex $beforeText
    if ($i -eq -1) {
        Write-Output "WARNING: could not locate paragraph: $beforeText"
        return
    }
    $d.Paragraphs.Item($i).Range.InsertParagraphBefore()
    # After InsertParagraphBefore, index $i is now the freshly-made blank
    # paragraph, and the original text has shifted down to index $i + 1.
    $d.Paragraphs.Item($i).Range.Text = $newText
}

# 1. Mick's first wordless reaction gets a mood tag.
Replace-Text "Mick: …" "Mick (arms_crossed skeptical): …"

# 2. Narration tweak: "he called for" -> "he wanted to talk to"
Replace-Text "He stares at me for a little longer, causing me to wonder if it really was me he called for." "He stares at me for a little longer, causing me to wonder if it really was me he wanted to talk to."

# 3. "How was practice?" gets a mood tag.
Replace-Text "Mick: How was practice?" "Mick (arms_crossed neutral): How was practice?"

# 4. "I see. That's good to hear." gets a mood tag.
Replace-Text "Mick: I see. That’s good to hear." "Mick (neutral indifferent): I see. That’s good to hear."

# 5. "Um…" gets a mood tag.
Replace-Text "Mick: Um…" "Mick (embarrassed neutral): Um…"

# 6. "Thanks. For going with her." gets a mood tag.
Replace-Text "Mick: Thanks. For going with her." "Mick (embarrassed embarrassed): Thanks. For going with her."

# 7. "Huh…?" becomes two beats: a standalone "…" paragraph, followed by "Huh?"
$huhIdx = Find-ParagraphIndex "Huh…?"
if ($huhIdx -eq -1) {
    Write-Output "WARNING: could not locate 'Huh...?' paragraph"
} else {
    $d.Paragraphs.Item($huhIdx).Range.InsertParagraphBefore()
    $d.Paragraphs.Item($huhIdx).Range.Text = "…"
    $d.Paragraphs.Item($huhIdx + 1).Range.Text = "Huh?"
}

# 8. "I quit playing an instrument..." gets a mood tag.
Replace-Text "Mick: I quit playing an instrument this year to join a club, so I feel a little bad about getting you to replace me." "Mick (embarrassed away): I quit playing an instrument this year to join a club, so I feel a little bad about getting you to replace me."

# 9. New standalone cue "Mick (neutral sigh):" right before Pro's "don't worry about it" line.
Insert-ParagraphBefore "Pro: Don’t worry about it. It’s not like I have anything better to do." "Mick (neutral sigh):"

# 10. New standalone cue "Mick (neutral neutral):" right before "He nods slowly..."
Insert-ParagraphBefore "He nods slowly, and I start to get the feeling that he’s as awkward as Prim." "Mick (neutral neutral):"

# 11. "Oh, right." gets a mood tag.
Replace-Text "Mick: Oh, right." "Mick (neutral curious): Oh, right."

# 12. "I'll see you around, I guess." gets a mood tag.
Replace-Text "Mick: I’ll see you around, I guess." "Mick (neutral neutral): I’ll see you around, I guess."

# 13. New standalone cue "Mick (exit):" right before "Mick ducks back into his classroom..."
Insert-ParagraphBefore "Mick ducks back into his classroom as I start to head back, hoping that Ms. Tran won’t give  me too much trouble if I show up late." "Mick (exit):"

Write-Output "done"
